$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.475.58'
$ws.Range('E2').Value = '  +0.71%  '
$ws.Range('D3').Value = '1.943.53'
$ws.Range('E3').Value = '  -2.00%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.82'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  -0.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.608'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  -3.46%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.66'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('E7').Value = '  -6.56%  '
$ws.Range('E9').Value = '  -2.92%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '55.92'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  -1.29%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0841'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').Value = '  +5.16%  '
$ws.Range('E12').Value = '  +0.44%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.829'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('E13').Value = '  -4.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.58'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  -4.18%  '
$ws.Range('D15').Value = '2.227.10'
$ws.Range('E15').Value = '  -2.06%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.65'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').Value = '  -2.91%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.25'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('E17').Value = '  -3.82%  '
$ws.Range('D18').Value = '1.966.45'
$ws.Range('E18').Value = '  -1.33%  '
$ws.Range('D19').Value = '36.385.71'
$ws.Range('E19').Value = '  +0.82%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.0₃0873'
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('B21').Value = 'Litecoin'
$ws.Range('C21').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '69.80'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  -2.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '230.06'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  -3.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.01'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  -5.15%  '
$ws.Range('E24').Value = '  +0.13%  '
$ws.Range('E25').Value = '  -1.62%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.30'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  -0.64%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.32'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  -4.64%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '162.68'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  +1.85%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.45'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  -2.19%  '
$ws.Range('E30').Value = '  -9.61%  '
$ws.Range('E31').Value = '  -1.61%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.16'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  +0.57%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.71'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  -4.49%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0634'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').Value = '  +1.55%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.29'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  -2.34%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.28'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  -1.14%  '
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('E38').Value = '  -3.18%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.16'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  -5.61%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.06'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  -2.44%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0974'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  -1.91%  '
$ws.Range('E42').Value = '  +4.30%  '
$ws.Range('E43').Value = '  -4.59%  '
$ws.Range('E44').Value = '  -2.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '16.11'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  -0.73%  '
$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.04'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  -5.53%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '1.354.14'
$ws.Range('E47').Value = '  -0.19%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '87.89'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  -5.85%  '
$ws.Range('E49').Value = '  -4.50%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.82'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  -1.68%  '
$ws.Range('E51').Value = '  +3.43%  '
